# people_widget.xlsx: change the "people" field from a multi-select widget to a
# plain textarea, and reformat its sample/default value to use a newline
# between names instead of a comma.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# C3 = type column for the "people" row -> select_or_add_multiple => textarea
$ws.Range("C3").Value = "textarea"

# M3 = value column for the "people" row -> "Ruben, Jessica" => "Ruben<LF>Jessica"
$nl = [char]10
$ws.Range("M3").Value = '"Ruben' + $nl + 'Jessica"'

# The new multi-line value needs wrap text turned on so it renders on the sheet
# (matches the style change from the non-wrapping style to the wrapping one).
$ws.Range("M3").WrapText = $true

# Move the active selection on the sheet from F3 to M3.
$ws.Range("M3").Select() | Out-Null
